$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell reference -> new display value, scraped from the updated coinranking.com snapshot.
$updates = @{
    'D2' = '36.964.58'
    'E2' = '  +0.05%  '
    'D3' = '2.081.10'
    'E3' = '  +8.58%  '
    'E4' = '  +0.01%  '
    'D5' = '251.25'
    'E5' = '  +0.57%  '
    'D6' = '0.658'
    'E6' = '  -6.03%  '
    'E7' = '  +0.06%  '
    'D8' = '50.71'
    'E8' = '  +5.76%  '
    'D9' = '60.56'
    'E9' = '  +4.32%  '
    'D10' = '0.376'
    'E10' = '  -0.61%  '
    'D11' = '0.0744'
    'E11' = '  -2.30%  '
    'D12' = '0.106'
    'E12' = '  +5.89%  '
    'D13' = '14.64'
    'E13' = '  -3.48%  '
    'D14' = '2.393.55'
    'E14' = '  +9.04%  '
    'D15' = '0.831'
    'E15' = '  +0.54%  '
    'D16' = '2.067.41'
    'E16' = '  +7.73%  '
    'D17' = '5.17'
    'E17' = '  +0.75%  '
    'D18' = '36.883.97'
    'E18' = '  -0.03%  '
    'D19' = '72.77'
    'E19' = '  -2.59%  '
    'D20' = '0.0₃0823'
    'E20' = '  -4.00%  '
    'D21' = '13.33'
    'E21' = '  -3.21%  '
    'D22' = '240.56'
    'E22' = '  -4.57%  '
    'D23' = '5.26'
    'E23' = '  +1.82%  '
    'E24' = '  -0.07%  '
    'D25' = '2.50'
    'E25' = '  +2.39%  '
    'D26' = '169.55'
    'E26' = '  +1.32%  '
    'D27' = '9.43'
    'E27' = '  +6.52%  '
    'D28' = '20.99'
    'E28' = '  +12.01%  '
    'D29' = '1.99'
    'E29' = '  -9.76%  '
    'B30' = 'Stellar'
    'C30' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D30' = '0.123'
    'E30' = '  -5.41%  '
    'B31' = 'Gas'
    'C31' = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
    'D31' = '24.62'
    'E31' = '  +26.55%  '
    'D32' = '1.15'
    'E32' = '  +30.17%  '
    'D33' = '4.50'
    'E33' = '  -2.30%  '
    'D34' = '0.0609'
    'E34' = '  -0.57%  '
    'D35' = '0.0915'
    'E35' = '  +2.24%  '
    'B36' = 'BinanceUSD'
    'C36' = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    'D36' = '1.00'
    'E36' = '  +0.03%  '
    'B37' = 'LidoDAOToken'
    'C37' = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    'D37' = '2.33'
    'E37' = '  +17.37%  '
    'D38' = '4.11'
    'E38' = '  -5.29%  '
    'D39' = '1.82'
    'E39' = '  -4.67%  '
    'D40' = '1.32'
    'E40' = '  -11.28%  '
    'D41' = '0.0225'
    'E41' = '  -1.51%  '
    'B42' = 'InjectiveProtocol'
    'C42' = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    'D42' = '17.37'
    'E42' = '  -2.19%  '
    'B43' = 'ARBITRUM'
    'C43' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D43' = '1.15'
    'E43' = '  +5.23%  '
    'B44' = 'Aave'
    'C44' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D44' = '97.38'
    'E44' = '  -7.32%  '
    'E45' = '  -2.96%  '
    'D46' = '1.344.78'
    'E46' = '  -1.32%  '
    'D47' = '0.0856'
    'E47' = '  +3.69%  '
    'D48' = '2.91'
    'E48' = '  +3.28%  '
    'D49' = '6.98'
    'E49' = '  +8.68%  '
    'B50' = 'MultiversX'
    'C50' = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
    'D50' = '47.12'
    'E50' = '  +9.27%  '
    'B51' = 'RocketPoolETH'
    'C51' = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
    'D51' = '2.265.90'
    'E51' = '  +8.14%  '
}

foreach ($ref in $updates.Keys) {
    $value = $updates[$ref]
    $cell = $ws.Range($ref)
    if ($ref.StartsWith("D") -and ($value -match '^-?\d+(\.\d+)?$')) {
        # Plain numeric-looking price strings must stay text (matches the sheet's
        # original inline-string cell type), so force text entry via a leading quote
        # and then drop back to the default "Normal" style (no custom number format).
        $cell.Value = "'" + $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
